$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "CODES": add a new error-code row (400 / GENERAL STARTUP ERROR) in
# the previously-blank row 33, matching the look of the other data rows
# (e.g. row 9: NUMBER CODE column style, NAME column style).
# ---------------------------------------------------------------------------
$codes = $wb.Worksheets.Item("CODES")

$codes.Range("A33").NumberFormat = "@"
$codes.Range("A33").Value = "400"
$codes.Range("A9").Copy()
$codes.Range("A33").PasteSpecial(-4122)

$codes.Range("B33").NumberFormat = "@"
$codes.Range("B33").Value = "GENERAL STARTUP ERROR"
$codes.Range("B9").Copy()
$codes.Range("B33").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Sheet "GENERAL RANGES DESC": add the two new general-range rows (4XX /
# STARTUP ERROR, and 5XX / RANDOM GENERATED ERRORS) right after the
# existing 3XX row, matching the style used by the other data rows.
# ---------------------------------------------------------------------------
$ranges = $wb.Worksheets.Item("GENERAL RANGES DESC")

$ranges.Range("A7").NumberFormat = "@"
$ranges.Range("A7").Value = "4XX"
$ranges.Range("A2").Copy()
$ranges.Range("A7").PasteSpecial(-4122)

$ranges.Range("B7").NumberFormat = "@"
$ranges.Range("B7").Value = "STARTUP ERROR"
$ranges.Range("B2").Copy()
$ranges.Range("B7").PasteSpecial(-4122)

$ranges.Range("A8").NumberFormat = "@"
$ranges.Range("A8").Value = "5XX"
$ranges.Range("A2").Copy()
$ranges.Range("A8").PasteSpecial(-4122)

$ranges.Range("B8").NumberFormat = "@"
$ranges.Range("B8").Value = "RANDOM GENERATED ERRORS"
$ranges.Range("B2").Copy()
$ranges.Range("B8").PasteSpecial(-4122)
